$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6933.3335
$ws.Range("I40").Value = 1300
$ws.Range("J40").Value = 9750
$ws.Range("K40").Value = 1300
$ws.Range("L40").Value = 9750
$ws.Range("M40").Value = -1125
$ws.Range("N40").Value = -10100

$ws.Range("H51").Value = 3059.8
$ws.Range("I51").Value = 1649.5
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 1649.5
$ws.Range("L51").Value = 4000
$ws.Range("M51").Value = -1165.5
$ws.Range("N51").Value = -4968

$ws.Range("H64").Value = 2984.6155
$ws.Range("I64").Value = 2900
$ws.Range("K64").Value = 2900
$ws.Range("M64").Value = -2652

$ws.Range("H67").Value = 2984.6155
$ws.Range("I67").Value = 2900
$ws.Range("K67").Value = 2900
$ws.Range("M67").Value = -2042

$ws.Range("H113").Value = 2800.5557
$ws.Range("I113").Value = 2641
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2641
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 613
$ws.Range("N113").Value = -9508

$ws.Range("H138").Value = 3147.7312
$ws.Range("I138").Value = 1513.1333
$ws.Range("J138").Value = 3462.077
$ws.Range("K138").Value = 4539.3999
$ws.Range("L138").Value = 10386.231
$ws.Range("M138").Value = 600.6000999999997
$ws.Range("N138").Value = -20666.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7590868.5
$ws.Range("I32").Value = 9104250
$ws.Range("K32").Value = 9104250
$ws.Range("M32").Value = -9103963

$ws.Range("H45").Value = 3688.5557
$ws.Range("I45").Value = 5506
$ws.Range("K45").Value = 5506
$ws.Range("M45").Value = -5129

$ws.Range("H93").Value = 64149.332
$ws.Range("J93").Value = 64149.332
$ws.Range("L93").Value = 64149.332
$ws.Range("N93").Value = -69141.33199999999

$ws.Range("H102").Value = 2625
$ws.Range("I102").Value = 2625
$ws.Range("K102").Value = 2625
$ws.Range("M102").Value = -1003

$ws.Range("H132").Value = 2657559.2
$ws.Range("I132").Value = 5918.5835
$ws.Range("J132").Value = 4529305.5
$ws.Range("K132").Value = 17755.7505
$ws.Range("L132").Value = 13587916.5
$ws.Range("M132").Value = -15225.7505
$ws.Range("N132").Value = -13592976.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 69021
$ws.Range("J21").Value = 69021
$ws.Range("L21").Value = 69021
$ws.Range("N21").Value = -69493

$ws.Range("H69").Value = 33443.89
$ws.Range("J69").Value = 33443.89
$ws.Range("L69").Value = 33443.89
$ws.Range("N69").Value = -35065.89

$ws.Range("H72").Value = 33443.89
$ws.Range("J72").Value = 33443.89
$ws.Range("L72").Value = 100331.67
$ws.Range("N72").Value = -108443.67

$ws.Range("H75").Value = 25207.334
$ws.Range("I75").Value = 3480
$ws.Range("J75").Value = 33564
$ws.Range("K75").Value = 3480
$ws.Range("L75").Value = 33564
$ws.Range("M75").Value = -2544
$ws.Range("N75").Value = -35436

$ws.Range("H78").Value = 25207.334
$ws.Range("I78").Value = 3480
$ws.Range("J78").Value = 33564
$ws.Range("K78").Value = 10440
$ws.Range("L78").Value = 100692
$ws.Range("M78").Value = -5760
$ws.Range("N78").Value = -110052

$ws.Range("H106").Value = 76557
$ws.Range("J106").Value = 76557
$ws.Range("L106").Value = 76557
$ws.Range("N106").Value = -79081

$ws.Range("H134").Value = 2984.1177
$ws.Range("I134").Value = 2357
$ws.Range("J134").Value = 3689.625
$ws.Range("K134").Value = 7071
$ws.Range("L134").Value = 11068.875
$ws.Range("M134").Value = -4536
$ws.Range("N134").Value = -16138.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5593.9854
$ws.Range("I31").Value = 1594.2693
$ws.Range("J31").Value = 8012.4185
$ws.Range("K31").Value = 1594.2693
$ws.Range("L31").Value = 8012.4185
$ws.Range("M31").Value = -1299.2693
$ws.Range("N31").Value = -8602.4185

$ws.Range("H34").Value = 5593.9854
$ws.Range("I34").Value = 1594.2693
$ws.Range("J34").Value = 8012.4185
$ws.Range("K34").Value = 1594.2693
$ws.Range("L34").Value = 8012.4185
$ws.Range("M34").Value = -1392.2693
$ws.Range("N34").Value = -8416.4185

$ws.Range("H62").Value = 4378.4077
$ws.Range("I62").Value = 4527.646
$ws.Range("K62").Value = 4527.646
$ws.Range("M62").Value = -3903.646

$ws.Range("H65").Value = 4378.4077
$ws.Range("I65").Value = 4527.646
$ws.Range("K65").Value = 22638.23
$ws.Range("M65").Value = -19518.23

$ws.Range("H81").Value = 97164
$ws.Range("J81").Value = 97164
$ws.Range("L81").Value = 97164
$ws.Range("N81").Value = -99160

$ws.Range("H84").Value = 97164
$ws.Range("J84").Value = 97164
$ws.Range("L84").Value = 291492
$ws.Range("N84").Value = -301476

$ws.Range("H99").Value = 2515.0967
$ws.Range("I99").Value = 2428.4285
$ws.Range("J99").Value = 2540.375
$ws.Range("K99").Value = 2428.4285
$ws.Range("L99").Value = 2540.375
$ws.Range("M99").Value = -930.4285
$ws.Range("N99").Value = -5536.375

$ws.Range("H107").Value = 3907038
$ws.Range("I107").Value = 5682282.5
$ws.Range("K107").Value = 5682282.5
$ws.Range("M107").Value = -5680362.5

$ws.Range("H126").Value = 2515.0967
$ws.Range("I126").Value = 2428.4285
$ws.Range("J126").Value = 2540.375
$ws.Range("K126").Value = 7285.2855
$ws.Range("L126").Value = 7621.125
$ws.Range("M126").Value = -4815.2855
$ws.Range("N126").Value = -12561.125

$ws.Range("H141").Value = 71267.87
$ws.Range("I141").Value = 107000
$ws.Range("J141").Value = 69282.75
$ws.Range("K141").Value = 107000
$ws.Range("L141").Value = 69282.75
$ws.Range("M141").Value = -101820
$ws.Range("N141").Value = -79642.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 4989
$ws.Range("J93").Value = 4994.154
$ws.Range("L93").Value = 14982.462
$ws.Range("N93").Value = -18726.462

$ws.Range("H107").Value = 2140.5789
$ws.Range("I107").Value = 419.7
$ws.Range("J107").Value = 4052.6667
$ws.Range("K107").Value = 1259.1
$ws.Range("L107").Value = 12158.0001
$ws.Range("M107").Value = 660.9000000000001
$ws.Range("N107").Value = -15998.0001

$ws.Range("H108").Value = 1663.5
$ws.Range("J108").Value = 2950
$ws.Range("L108").Value = 8850
$ws.Range("N108").Value = -14610

$ws.Range("H113").Value = 561.4074000000001
$ws.Range("I113").Value = 583.8461
$ws.Range("J113").Value = 540.5714
$ws.Range("K113").Value = 1751.5383
$ws.Range("L113").Value = 1621.7142
$ws.Range("M113").Value = 418.4617000000001
$ws.Range("N113").Value = -5961.7142

$ws.Range("H134").Value = 6162.404
$ws.Range("I134").Value = 2380.9048
$ws.Range("K134").Value = 7142.714399999999
$ws.Range("M134").Value = -2072.714399999999

$ws.Range("H136").Value = 2390.6428
$ws.Range("I136").Value = 847.4167
$ws.Range("J136").Value = 11650
$ws.Range("K136").Value = 2542.2501
$ws.Range("L136").Value = 34950
$ws.Range("M136").Value = 2557.7499
$ws.Range("N136").Value = -45150

$ws.Range("H137").Value = 34024.17
$ws.Range("J137").Value = 56331.105
$ws.Range("L137").Value = 168993.315
$ws.Range("N137").Value = -179193.315

$ws.Range("H139").Value = 225479.16
$ws.Range("I139").Value = 358738.94
$ws.Range("J139").Value = 5992.4707
$ws.Range("K139").Value = 1076216.82
$ws.Range("L139").Value = 17977.4121
$ws.Range("M139").Value = -1071076.82
$ws.Range("N139").Value = -28257.4121

$ws.Range("H140").Value = 1763.4412
$ws.Range("I140").Value = 1160.85
$ws.Range("K140").Value = 3482.55
$ws.Range("M140").Value = 1697.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 90357.164
$ws.Range("J86").Value = 90357.164
$ws.Range("L86").Value = 90357.164
$ws.Range("N86").Value = -92729.164

$ws.Range("H89").Value = 90357.164
$ws.Range("J89").Value = 90357.164
$ws.Range("L89").Value = 271071.492
$ws.Range("N89").Value = -282927.492

$ws.Range("H132").Value = 58832900
$ws.Range("I132").Value = 100012530
$ws.Range("J132").Value = 4860.7144
$ws.Range("K132").Value = 300037590
$ws.Range("L132").Value = 14582.1432
$ws.Range("M132").Value = -300035060
$ws.Range("N132").Value = -19642.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1929
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2000.5
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2000.5
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -3498.5

$ws.Range("H71").Value = 1929
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2000.5
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 10002.5
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -17490.5

$ws.Range("H132").Value = 2899.818
$ws.Range("I132").Value = 2605.8235
$ws.Range("K132").Value = 7817.470499999999
$ws.Range("M132").Value = -5287.470499999999

$ws.Range("H136").Value = 2105.75
$ws.Range("I136").Value = 1958.7
$ws.Range("J136").Value = 2350.8333
$ws.Range("K136").Value = 5876.1
$ws.Range("L136").Value = 7052.499899999999
$ws.Range("M136").Value = -3326.1
$ws.Range("N136").Value = -12152.4999

$ws.Range("H140").Value = 57224.418
$ws.Range("J140").Value = 57224.418
$ws.Range("L140").Value = 57224.418
$ws.Range("N140").Value = -67584.41800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9725710
$ws.Range("I132").Value = 3770.6316
$ws.Range("J132").Value = 26518152
$ws.Range("K132").Value = 11311.8948
$ws.Range("L132").Value = 79554456
$ws.Range("M132").Value = -8781.8948
$ws.Range("N132").Value = -79559516

$ws.Range("H136").Value = 3631.475
$ws.Range("I136").Value = 3475
$ws.Range("J136").Value = 3892.2666
$ws.Range("K136").Value = 10425
$ws.Range("L136").Value = 11676.7998
$ws.Range("M136").Value = -7875
$ws.Range("N136").Value = -16776.7998
